$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3628304004669189
$ws.Range("B1").Value = 2.827953100204468
$ws.Range("C1").Value = 4.700586795806885
$ws.Range("D1").Value = 1.736625194549561
$ws.Range("E1").Value = 0.8260374665260315
